$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.409.39"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "1.849.84"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.72"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6302"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07670"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.76%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2941"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.49"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07750"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").Value = "1.839.74"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.019"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.54%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001090"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +8.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6802"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.61"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.95%  "

$ws.Range("D17").Value = "2.095.12"
$ws.Range("E17").Value = "  -7.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.150"
$ws.Range("D18").ClearFormats()

$ws.Range("D19").Value = "29.423.46"
$ws.Range("E19").Value = "  +0.08%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.33"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.70%  "

$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.443"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.28%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.51"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.364"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("E28").Value = "  +0.26%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.468"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.312"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05684"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.113"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.050"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.159"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.38%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7097"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.588"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.781"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.03%  "

$ws.Range("D39").Value = "1.232.40"
$ws.Range("E39").Value = "  -2.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01797"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.475"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9143"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.49"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.43%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.21"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000121"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.21%  "

$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.164"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.81%  "

$ws.Range("B48").Value = "TheSandbox"
$ws.Range("C48").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4015"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.16%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.057"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.689"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.09%  "

$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1126"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.11%  "
